$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.425.66"
$ws.Range("E2").Value = "  +2.88%  "
$ws.Range("D3").Value = "2.015.74"
$ws.Range("E3").Value = "  +6.30%  "
$ws.Range("E4").Value = "  -0.02%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "245.46"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.47%  "
$ws.Range("E6").Value = "  -4.62%  "
$ws.Range("E7").Value = "  +0.01%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "45.05"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +4.56%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "59.96"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +6.73%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.369"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +3.26%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0713"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -5.37%  "
$ws.Range("E12").Value = "  +0.38%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "14.55"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +3.05%  "
$ws.Range("D14").Value = "2.312.41"
$ws.Range("E14").Value = "  +6.43%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.810"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +1.53%  "
$ws.Range("D16").Value = "2.020.05"
$ws.Range("E16").Value = "  +7.20%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "4.90"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -2.10%  "
$ws.Range("D18").Value = "36.262.95"
$ws.Range("E18").Value = "  +2.22%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "71.29"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -3.23%  "
$ws.Range("D20").Value = "0.0₃0816"
$ws.Range("E20").Value = "  -1.79%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "12.90"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -0.63%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "236.12"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -3.42%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "4.87"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -6.42%  "
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("E25").Value = "  -8.71%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "163.68"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -2.18%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "19.63"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +7.14%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "8.58"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -0.68%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "1.93"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -10.94%  "
$ws.Range("E30").Value = "  -4.67%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "22.58"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +64.39%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "4.41"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +1.28%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.0589"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -2.36%  "
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("E35").Value = "  -0.74%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "4.00"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -5.70%  "
$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "2.16"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +10.49%  "
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.0804"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +8.65%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.849"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -0.82%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "1.33"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -9.50%  "
$ws.Range("E41").Value = "  -3.83%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "95.99"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -3.22%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "1.12"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +2.99%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "2.78"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +14.78%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "15.95"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -6.13%  "
$ws.Range("D46").Value = "1.318.35"
$ws.Range("E46").Value = "  -0.63%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.0818"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +0.96%  "
$ws.Range("E48").Value = "  +1.10%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "2.21"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -6.49%  "
$ws.Range("D50").Value = "2.202.58"
$ws.Range("E50").Value = "  +6.18%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "3.86"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +15.66%  "
